$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.220.89"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").Value = "1.650.78"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.45"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.516"
$ws.Range("E6").Value = "  +2.23%  "
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("E8").Value = "  +1.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0630"
$ws.Range("E9").Value = "  +1.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.04"
$ws.Range("E10").Value = "  +1.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0849"
$ws.Range("E11").Value = "  +0.51%  "
$ws.Range("D12").Value = "1.880.99"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").Value = "1.644.35"
$ws.Range("E13").Value = "  -0.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.16"
$ws.Range("E14").Value = "  +0.39%  "
$ws.Range("E15").Value = "  +2.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.73"
$ws.Range("E16").Value = "  +1.82%  "
$ws.Range("D17").Value = "27.211.71"
$ws.Range("E17").Value = "  +1.04%  "
$ws.Range("D18").Value = "0.0₃0740"
$ws.Range("E18").Value = "  +1.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "220.32"
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.86"
$ws.Range("E21").Value = "  +2.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.58"
$ws.Range("E22").Value = "  +6.31%  "
$ws.Range("E23").Value = "  +0.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.24"
$ws.Range("E24").Value = "  +0.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.32"
$ws.Range("E25").Value = "  +0.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.55"
$ws.Range("E26").Value = "  +2.62%  "
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.82"
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0510"
$ws.Range("E30").Value = "  -0.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.19"
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.41"
$ws.Range("E32").Value = "  +0.59%  "
$ws.Range("E33").Value = "  +1.49%  "
$ws.Range("E34").Value = "  +1.71%  "
$ws.Range("D35").Value = "1.265.99"
$ws.Range("E35").Value = "  +1.24%  "
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("E37").Value = "  +1.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.549"
$ws.Range("E38").Value = "  +3.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.847"
$ws.Range("E39").Value = "  +1.83%  "
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.811"
$ws.Range("E41").Value = "  +0.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.45"
$ws.Range("E42").Value = "  +1.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.23"
$ws.Range("E43").Value = "  +5.66%  "
$ws.Range("D44").Value = "1.791.47"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.15"
$ws.Range("E45").Value = "  +1.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.75"
$ws.Range("E46").Value = "  +0.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.61"
$ws.Range("E47").Value = "  +0.61%  "
$ws.Range("D48").Value = "0.0₆0105"
$ws.Range("E48").Value = "  +0.82%  "
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("E50").Value = "  +0.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0975"
$ws.Range("E51").Value = "  +0.13%  "
